$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 65, shifting existing rows 65:133 down to 66:134
$ws.Rows("65:65").Insert()

# Populate the newly inserted row 65 with the new data record.
# Columns A, B, C, E, F, G, H, I, N, O, Q, R are identical across this block,
# so copy them from row 66 (the row that used to be row 65 before the insert).
$ws.Range("A65").Value = $ws.Range("A66").Value2
$ws.Range("B65").Value = $ws.Range("B66").Value2
$ws.Range("C65").Value = $ws.Range("C66").Value2
$ws.Range("D65").Value = 44658
$ws.Range("E65").Value = $ws.Range("E66").Value2
$ws.Range("F65").Value = $ws.Range("F66").Value2
$ws.Range("G65").Value = $ws.Range("G66").Value2
$ws.Range("H65").Value = $ws.Range("H66").Value2
$ws.Range("I65").Value = $ws.Range("I66").Value2
$ws.Range("J65").Value = 60
$ws.Range("K65").Value = 9000
$ws.Range("L65").Value = 9000
$ws.Range("M65").Value = 9000
$ws.Range("N65").Value = $ws.Range("N66").Value2
$ws.Range("O65").Value = $ws.Range("O66").Value2
$ws.Range("P65").Value = 3000
$ws.Range("Q65").Value = $ws.Range("Q66").Value2
$ws.Range("R65").Value = $ws.Range("R66").Value2
